$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Grayson Allen", "PG,SG,SF", "Phoenix Suns"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Wendell Carter Jr.", "PF,C", "Orlando Magic"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Mike Conley", "PG", "Minnesota Timberwolves")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
